$wb = $excel.ActiveWorkbook

# --- capex sheet: row 2 (BF-BOF) values 50 -> 30 across B2:AA2 ---
$wsCapex = $wb.Worksheets.Item("capex")
$wsCapex.Range("B2:AA2").Value = 30
$wsCapex.Activate()
$wsCapex.Range("O17").Select() | Out-Null

# --- baseline sheet: C4 fuel changed from "Electricity" to "Bio" ---
$wsBaseline = $wb.Worksheets.Item("baseline")
$wsBaseline.Range("C4").Value = "Bio"
$wsBaseline.Activate()
$wsBaseline.Range("F2").Select() | Out-Null

# --- technology sheet: lifespan/introduced_year updates ---
$wsTech = $wb.Worksheets.Item("technology")
$wsTech.Range("B2").Value = 20
$wsTech.Range("C2").Value = 2020
$wsTech.Range("B4").Value = 20
$wsTech.Range("C4").Value = 2020

# --- fuel_cost sheet: just a selection change ---
$wsFuelCost = $wb.Worksheets.Item("fuel_cost")
$wsFuelCost.Activate()
$wsFuelCost.Range("B2").Select() | Out-Null

# Make technology sheet the active / selected tab (mirrors tabSelected move)
$wsTech.Activate()
$wsTech.Range("B5").Select() | Out-Null
